$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha, serial date), J (Volumen),
# K (Precio minimo), L (Precio maximo), M (Precio promedio ponderado),
# P (Precio $/Kg). This reflects a reshuffle of the weekly rows.

$data = @{
    2  = @{ D = 44915; J = 50; K = 18000; L = 18000; M = 18000; P = 1385 }
    3  = @{ D = 44832; J = 60; K = 17000; L = 18000; M = 17500; P = 1346 }
    4  = @{ D = 44930; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    5  = @{ D = 45155; J = 30; K = 20000; L = 20000; M = 20000; P = 1538 }
    6  = @{ D = 44839; J = 40; K = 15000; L = 16000; M = 15500; P = 1192 }
    7  = @{ D = 44868; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    10 = @{ D = 44880; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    11 = @{ D = 44859; J = 30; K = 13000; L = 13000; M = 13000; P = 1000 }
    12 = @{ D = 44841; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    13 = @{ D = 44894; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    14 = @{ D = 44874; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    15 = @{ D = 44895; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    16 = @{ D = 44922; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    17 = @{ D = 44943; J = 30; K = 17000; L = 17000; M = 17000; P = 1308 }
    18 = @{ D = 44846; J = 30; K = 18000; L = 18000; M = 18000; P = 1385 }
    19 = @{ D = 44797; J = 60; K = 12000; L = 13000; M = 12500; P = 962 }
    20 = @{ D = 44810; J = 40; K = 12000; L = 13000; M = 12500; P = 962 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value2 = $vals.D
    $ws.Cells.Item($row, 10).Value2 = $vals.J
    $ws.Cells.Item($row, 11).Value2 = $vals.K
    $ws.Cells.Item($row, 12).Value2 = $vals.L
    $ws.Cells.Item($row, 13).Value2 = $vals.M
    $ws.Cells.Item($row, 16).Value2 = $vals.P
}
